# Weekly update: insert a new price-report pair (Primera / Segunda) at the
# top of the historical data block for Betarraga - Terminal La Palmera de
# La Serena, shifting the rest of the rows down by two and pushing the two
# oldest rows out to the end of the table (rows 196/197).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new pair of rows by inserting two blank rows at the
# current position of the oldest-but-one record (row 54). Excel shifts
# rows 54:195 down to 56:197 and extends the used range accordingly.
$ws.Range("A54:R55").Insert()

# New "Primera" quality record.
$ws.Range("A54").Value = 8
$ws.Range("B54").Value = "Terminal La Palmera de La Serena"
$ws.Range("C54").Value = "Coquimbo"
$ws.Range("D54").Value = 44525
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100114014
$ws.Range("G54").Value = "Betarraga"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 3000
$ws.Range("K54").Value = 450
$ws.Range("L54").Value = 500
$ws.Range("M54").Value = 475
$ws.Range("N54").Value = "`$/paquete 3 unidades"
$ws.Range("O54").Value = "Provincia del Elquí"
$ws.Range("P54").Value = 158
$ws.Range("Q54").Value = 3
$ws.Range("R54").Value = "Hortaliza"

# New "Segunda" quality record.
$ws.Range("A55").Value = 8
$ws.Range("B55").Value = "Terminal La Palmera de La Serena"
$ws.Range("C55").Value = "Coquimbo"
$ws.Range("D55").Value = 44525
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 100114014
$ws.Range("G55").Value = "Betarraga"
$ws.Range("H55").Value = "Sin especificar"
$ws.Range("I55").Value = "Segunda"
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 350
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = 375
$ws.Range("N55").Value = "`$/paquete 3 unidades"
$ws.Range("O55").Value = "Provincia del Elquí"
$ws.Range("P55").Value = 125
$ws.Range("Q55").Value = 3
$ws.Range("R55").Value = "Hortaliza"
